# "data provider added in test" - update the Tweets test-data sheet
# with richer sample strings, used by the data-provider driven test.

$wb = $excel.ActiveWorkbook
$wsTweets = $wb.Worksheets.Item("Sheet2")

# Make sure we are working on / looking at the Tweets sheet (it is already
# the active tab in the workbook).
$wsTweets.Activate()

# Update the shared-string values used by the data provider.
$wsTweets.Range("A2").Value = "Hello All Indians"
$wsTweets.Range("A3").Value = "Good Evening India and Noida"
$wsTweets.Range("A4").Value = "Automation Framework TDD"

# Ensure the fonts used throughout the workbook carry an explicit
# (default) charset, matching how Excel normally persists font info.
$wb.Worksheets.Item("Sheet1").Cells.Font.Charset = 1
$wsTweets.Cells.Font.Charset = 1

# The selection on the Tweets sheet moved from A5 to A4.
$wsTweets.Range("A4").Select()

# Column A is now wider to fit the longer sample text.
$wsTweets.Columns.Item(1).ColumnWidth = 25.15

# First page number is no longer forced to start at 1.
$wsTweets.PageSetup.FirstPageNumber = 0
